$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.339398838185446
$ws.Range("C2").Value = -2.325962872140053
$ws.Range("D2").Value = 2.717184218186617
$ws.Range("E2").Value = -2.490917476823483
$ws.Range("F2").Value = 5.433053069595729
$ws.Range("B3").Value = 2.354781502736316
$ws.Range("C3").Value = -4.078462363779821
$ws.Range("D3").Value = 4.661151271876534
$ws.Range("E3").Value = 6.651705714458217
$ws.Range("F3").Value = 6.709178096077707
$ws.Range("G3").Value = -3.888207370780819
$ws.Range("H3").Value = 2.709991627146224
$ws.Range("B4").Value = 1.389778685219852
$ws.Range("C4").Value = 4.554594940432864
$ws.Range("D4").Value = 5.19183467676153
$ws.Range("E4").Value = -4.830094033492927
$ws.Range("F4").Value = 1.941824996538614
$ws.Range("B5").Value = 2.8698272557899
$ws.Range("C5").Value = -4.494149346636093
$ws.Range("D5").Value = 2.156850403891937
$ws.Range("E5").Value = -5.242968891207937
$ws.Range("F5").Value = -6.916046429188043
$ws.Range("G5").Value = -1.058669429628892
$ws.Range("H5").Value = 0.08183550252746841
$ws.Range("B6").Value = 4.6132376685254
$ws.Range("C6").Value = -4.323059390703492
$ws.Range("D6").Value = -6.47192035718534
$ws.Range("E6").Value = -1.133642100495023
$ws.Range("F6").Value = -0.1967381681158997
$ws.Range("B7").Value = -5.97921302499303
$ws.Range("C7").Value = -0.3336421004950549
$ws.Range("D7").Value = 0.5032618318840889
$ws.Range("E7").Value = -0.1632848153601001
$ws.Range("F7").Value = -1.291754278273445
$ws.Range("G7").Value = -0.8433771880326331
$ws.Range("H7").Value = 2.58809437072307
$ws.Range("B8").Value = 1.187206647323194
$ws.Range("C8").Value = -0.4360795331101601
$ws.Range("D8").Value = -1.691754278273436
$ws.Range("E8").Value = -1.143377188032644
$ws.Range("F8").Value = 2.388094370723039
$ws.Range("B9").Value = -0.9993369614038421
$ws.Range("C9").Value = -0.9162886220676389
$ws.Range("D9").Value = 2.379009745974372
$ws.Range("E9").Value = 1.681290410009282
$ws.Range("F9").Value = -0.3607148982294461
$ws.Range("G9").Value = 0.922161731270734
$ws.Range("H9").Value = 3.318707023947993
$ws.Range("B10").Value = 1.18809437072305
$ws.Range("C10").Value = 0.5777282364878857
$ws.Range("D10").Value = -1.038969365349572
$ws.Range("E10").Value = 0.456124862885531
$ws.Range("F10").Value = 2.812661112537114
$ws.Range("B11").Value = -0.943724682173837
$ws.Range("C11").Value = -0.5438751371144408
$ws.Range("D11").Value = 2.312661112537128
$ws.Range("E11").Value = 3.700005087957896
$ws.Range("F11").Value = -0.3330793704473369
$ws.Range("G11").Value = 0.4882750851338071
$ws.Range("H11").Value = 0.01467091632704909
$ws.Range("B12").Value = 1.211916711296126
$ws.Range("C12").Value = 2.499975276622385
$ws.Range("D12").Value = -0.8330762883444149
$ws.Range("E12").Value = 0.6882924312177612
$ws.Range("F12").Value = 0.3146750467264641
$ws.Range("B13").Value = 0.0669210872980841
$ws.Range("C13").Value = 0.2882771288503023
$ws.Range("D13").Value = 0.1146719783196859
$ws.Range("E13").Value = -0.7085072997627999
$ws.Range("F13").Value = 0.4807953866492261
$ws.Range("G13").Value = 0.02220952249070907
$ws.Range("H13").Value = -0.05749613082447191
$ws.Range("B14").Value = 0.007307160242419952
$ws.Range("C14").Value = -0.526487056322295
$ws.Range("D14").Value = 0.6661963638293009
$ws.Range("E14").Value = 0.1877297131885172
$ws.Range("F14").Value = 0.09679325291892904
$ws.Range("B15").Value = -0.1067066338104181
$ws.Range("C15").Value = -0.04323335916021187
$ws.Range("D15").Value = -0.166605428577582
$ws.Range("E15").Value = 0.2438558147322278
$ws.Range("F15").Value = 1.21360082768301
$ws.Range("G15").Value = 0.6139174841733281
$ws.Range("H15").Value = 1.240880967708379
$ws.Range("B16").Value = -0.9603412467216681
$ws.Range("C16").Value = -0.3615542135854211
$ws.Range("D16").Value = 0.9630191732715327
$ws.Range("E16").Value = 0.5052801541909131
$ws.Range("F16").Value = 1.229815784184155
$ws.Range("B17").Value = 0.5114728330376579
$ws.Range("C17").Value = 0.350496246476965
$ws.Range("D17").Value = 1.005295191856816
$ws.Range("E17").Value = 0.03764606290536698
$ws.Range("F17").Value = 1.298055291346585
$ws.Range("G17").Value = 1.13685241246894
$ws.Range("H17").Value = -1.191311050950602
$ws.Range("B18").Value = 0.8251185114959585
$ws.Range("C18").Value = 0.03207815415970594
$ws.Range("D18").Value = 1.29218462318174
$ws.Range("E18").Value = 1.162246220662766
$ws.Range("F18").Value = -1.216012244744192
$ws.Range("G18").Value = 0.8458859275583706
$ws.Range("H18").Value = -1.183478699330493
$ws.Range("I18").Value = 0.2823547203913674
$ws.Range("J18").Value = -0.8692897857011141
$ws.Range("B19").Value = 0.6426636127604581
$ws.Range("C19").Value = 1.138869314174826
$ws.Range("D19").Value = -1.189390679284543
$ws.Range("E19").Value = 0.8021181534470256
$ws.Range("F19").Value = -1.196748302637431
$ws.Range("G19").Value = 0.3079620948331154
$ws.Range("H19").Value = -0.7892897857011301
$ws.Range("B20").Value = -1.255996022051292
$ws.Range("C20").Value = 0.8995643220687695
$ws.Range("D20").Value = -1.20835213404358
$ws.Range("E20").Value = 0.3079620948330584
$ws.Range("F20").Value = -0.789289785701087
$ws.Range("G20").Value = 2.295090611800802
$ws.Range("H20").Value = -0.5099739529217828
$ws.Range("I20").Value = -0.09999999999999698
$ws.Range("J20").Value = 0.4654042787202852
$ws.Range("B21").Value = -0.7594837991400569
$ws.Range("C21").Value = 0.4399412104564964
$ws.Range("D21").Value = -0.734533228560807
$ws.Range("E21").Value = 2.275556695716489
$ws.Range("F21").Value = -0.5320054475438669
$ws.Range("G21").Value = -0.1220145169862861
$ws.Range("H21").Value = 0.4434067299780222
$ws.Range("B22").Value = -0.19336818743102
$ws.Range("C22").Value = 2.354131640238677
$ws.Range("D22").Value = -0.461748558094911
$ws.Range("E22").Value = -0.06732255871831994
$ws.Range("F22").Value = 0.4825157701251892
$ws.Range("G22").Value = 0.517677114878327
$ws.Range("H22").Value = 1.491403243014375
$ws.Range("I22").Value = 1.148757636989984
$ws.Range("J22").Value = -0.0544928672541829
$ws.Range("B23").Value = 2.422210334885136
$ws.Range("C23").Value = -0.5760769225611289
$ws.Range("D23").Value = -0.1439667268341251
$ws.Range("E23").Value = 0.4215053706022762
$ws.Range("F23").Value = 0.4562484426560931
$ws.Range("G23").Value = 1.430281299027996
$ws.Range("H23").Value = 1.087665430760481
$ws.Range("I23").Value = -0.1155553792681729
$ws.Range("B24").Value = -0.5511376795560359
$ws.Range("C24").Value = -0.22184360803742
$ws.Range("D24").Value = 0.4217508427945231
$ws.Range("E24").Value = 0.4954046212014191
$ws.Range("F24").Value = 1.469580582623763
$ws.Range("G24").Value = 1.127002351592182
$ws.Range("H24").Value = -0.07618089288315888
$ws.Range("B25").Value = -0.2239417963794241
$ws.Range("C25").Value = 0.5010313511557742
$ws.Range("D25").Value = 0.4970082786298951
$ws.Range("E25").Value = 1.470782412094053
$ws.Range("F25").Value = 1.028202489427699
$ws.Range("G25").Value = -0.1749823445367018
$ws.Range("H25").Value = 2.496480684450609
$ws.Range("I25").Value = 16.92263857113068
$ws.Range("B26").Value = -0.4345957212796748
$ws.Range("C26").Value = 0.5
$ws.Range("D26").Value = 1.474014352124283
$ws.Range("E26").Value = 1.031330924713933
$ws.Range("F26").Value = -0.1719573692094458
$ws.Range("G26").Value = 2.5
$ws.Range("H26").Value = 16.92569689670956
$ws.Range("B27").Value = 1.70590956953757
$ws.Range("C27").Value = 1.430893706027206
$ws.Range("D27").Value = 1.02751370927065
$ws.Range("E27").Value = -0.175670149342793
$ws.Range("F27").Value = 2.494083200150362
$ws.Range("G27").Value = 16.92124697637898
$ws.Range("B28").Value = -0.1259856478757172
$ws.Range("C28").Value = 1.031330924713933
$ws.Range("D28").Value = -0.1719573692094458
$ws.Range("E28").Value = 2.5
$ws.Range("F28").Value = 16.92569689670956
$ws.Range("B29").Value = 0.8276572911504871
$ws.Range("C29").Value = -0.217218108648638
$ws.Range("D29").Value = 2.538632612742652
$ws.Range("E29").Value = 16.96418152364336
$ws.Range("B30").Value = -0.2165566629339399
$ws.Range("C30").Value = 2.654937875620226
$ws.Range("D30").Value = 16.96416332357788
$ws.Range("B31").Value = 1.169600765142468
$ws.Range("C31").Value = 14.75546877091031
$ws.Range("B32").Value = 1.525696896709562
